# Contract_Progression_2425.xlsx -- "Adjustment et progression 12/12/24"
#
# Updates progress percentages and comment text on Sheet1 to reflect the
# 12/12/24 status update, and moves the on-screen selection down to G11
# (scrolled so row 6 is at the top), matching the author's final view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5: core HANWASH stakeholders / mWater account naming task ---
# status 73% -> 96%
$ws.Range("G5").Value = 0.96

# comment: "12 accounts from 43 need to be created to complete this task"
#       -> "2 accounts from 48 need to be created to complete this task"
# (keep the leading number/range phrase bold, rest of the sentence regular,
# matching the existing rich-text pattern used in that cell)
$h5 = $ws.Range("H5")
$h5.Value = "2 accounts from 48 need to be created to complete this task"
$h5.Font.Bold = $false
$h5.Characters(1, 18).Font.Bold = $true

# --- Row 8: Data quality check task ---
# status 25% -> 50%
$ws.Range("G8").Value = 0.5

# --- Row 10: User guides for the CPE console task ---
# status 20% -> 60%
$ws.Range("G10").Value = 0.6

# comment: "user guides is completed, not validated yet" -> "user guides is completed"
$ws.Range("H10").Value = "user guides is completed"

# --- Row 11: Training and support task ---
# status 50% -> 60%
$ws.Range("G11").Value = 0.6

# --- View state: scroll down a bit and move the selection to G11 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("G11").Select()
